# Remove column from alcohol measurement data (Sheet1, column M / index 13).
# Deleting the column shifts the former column N left into M, which is
# exactly the cell-value change captured by the diff (row by row, M<-N,
# N removed) and updates the sheet dimension from A1:N119 to A1:M119.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Columns.Item(13).Delete() | Out-Null

# The resave that produced this commit also renormalised the column widths
# (to a single, slightly narrower uniform width) on every sheet. Re-create
# that as closely as the object model allows by setting the width of the
# columns that actually carry formatting/data.
$ws1.Range($ws1.Cells.Item(1,1), $ws1.Cells.Item(1,13)).EntireColumn.ColumnWidth = 10.51

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Columns.Item(1).ColumnWidth = 10.51

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Columns.Item(1).ColumnWidth = 10.51

# Same resave also dropped the zoom level on every sheet's view from 65%/100%
# down to a shared 95%, and moved Sheet1's selection from AC1 to M1 (the new
# last used column after the delete). Gridlines stay visible throughout, so
# re-assert that explicitly too.
foreach ($name in @("Sheet1", "Sheet2", "Sheet3")) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Activate()
    $excel.ActiveWindow.Zoom = 95
    $excel.ActiveWindow.DisplayGridlines = $true
}

$ws1.Activate()
$ws1.Range("M1").Select() | Out-Null
